# Updates cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values below are plain decimal numbers (e.g. "25.12"); force them
# to remain text cells (matching the original inline-string storage) instead of
# letting Excel auto-convert them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = "70.451.51"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "2.520.46"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "577.11"
$ws.Range("E5").Value = "  -3.30%  "
$ws.Range("D6").Value = "168.93"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "2.520.42"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "2.983.95"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").Value = "70.418.73"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "25.12"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "2.521.09"
$ws.Range("E18").Value = "  -4.80%  "
$ws.Range("D19").Value = "11.45"
$ws.Range("E19").Value = "  -7.60%  "
$ws.Range("E20").Value = "  -6.34%  "
$ws.Range("D21").Value = "359.25"
$ws.Range("E21").Value = "  -3.05%  "
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("D23").Value = "1.98"
$ws.Range("E23").Value = "  -5.72%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "69.41"
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("D26").Value = "4.06"
$ws.Range("E26").Value = "  -6.14%  "
$ws.Range("D27").Value = "9.13"
$ws.Range("E27").Value = "  -7.11%  "
$ws.Range("D28").Value = "2.652.34"
$ws.Range("E28").Value = "  -4.91%  "
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "0.0₃0916"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("D32").Value = "485.04"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "156.00"
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "18.65"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "18.90"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").Value = "0.321"
$ws.Range("E42").Value = "  -3.59%  "
$ws.Range("E43").Value = "  -6.66%  "
$ws.Range("E44").Value = "  -12.74%  "
$ws.Range("E45").Value = "  -7.66%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").Value = "143.83"
$ws.Range("E47").Value = "  -8.16%  "
$ws.Range("D48").Value = "3.54"
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("D49").Value = "0.528"
$ws.Range("E49").Value = "  -5.53%  "
$ws.Range("E50").Value = "  -6.48%  "
$ws.Range("D51").Value = "0.599"
$ws.Range("E51").Value = "  -1.17%  "

